# Adds a "metadata" worksheet (sheetId 2) after the existing "data" sheet,
# refreshes the "panel query" timestamps recorded in column F of "data"
# (re-run of the panelapp export a bit later the same day), and populates
# the new "metadata" sheet with the panel-level summary row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the per-gene "time_taken" timestamps on the "data" sheet.
#    (same export re-run, ~3h43m later, same tight per-row spacing)
# ---------------------------------------------------------------------
$ts = @(
  "2021-10-05 14:35:49.056492",
  "2021-10-05 14:35:49.056499",
  "2021-10-05 14:35:49.056503",
  "2021-10-05 14:35:49.056505",
  "2021-10-05 14:35:49.056508",
  "2021-10-05 14:35:49.056511",
  "2021-10-05 14:35:49.056513",
  "2021-10-05 14:35:49.056516",
  "2021-10-05 14:35:49.056518",
  "2021-10-05 14:35:49.056521",
  "2021-10-05 14:35:49.056523",
  "2021-10-05 14:35:49.056526",
  "2021-10-05 14:35:49.056528",
  "2021-10-05 14:35:49.056531",
  "2021-10-05 14:35:49.056533",
  "2021-10-05 14:35:49.056535",
  "2021-10-05 14:35:49.056538",
  "2021-10-05 14:35:49.056541",
  "2021-10-05 14:35:49.056543",
  "2021-10-05 14:35:49.056546",
  "2021-10-05 14:35:49.056548",
  "2021-10-05 14:35:49.056551",
  "2021-10-05 14:35:49.056553",
  "2021-10-05 14:35:49.056555",
  "2021-10-05 14:35:49.056558",
  "2021-10-05 14:35:49.056561",
  "2021-10-05 14:35:49.056563",
  "2021-10-05 14:35:49.056566",
  "2021-10-05 14:35:49.056568",
  "2021-10-05 14:35:49.056571",
  "2021-10-05 14:35:49.056573",
  "2021-10-05 14:35:49.056575",
  "2021-10-05 14:35:49.056578",
  "2021-10-05 14:35:49.056581",
  "2021-10-05 14:35:49.056583",
  "2021-10-05 14:35:49.056586",
  "2021-10-05 14:35:49.056588",
  "2021-10-05 14:35:49.056591",
  "2021-10-05 14:35:49.056594",
  "2021-10-05 14:35:49.056596",
  "2021-10-05 14:35:49.056599",
  "2021-10-05 14:35:49.056601",
  "2021-10-05 14:35:49.056604",
  "2021-10-05 14:35:49.056606",
  "2021-10-05 14:35:49.056609",
  "2021-10-05 14:35:49.056611",
  "2021-10-05 14:35:49.056614",
  "2021-10-05 14:35:49.056616",
  "2021-10-05 14:35:49.056619",
  "2021-10-05 14:35:49.056621",
  "2021-10-05 14:35:49.056623",
  "2021-10-05 14:35:49.056626",
  "2021-10-05 14:35:49.056629",
  "2021-10-05 14:35:49.056631",
  "2021-10-05 14:35:49.056634",
  "2021-10-05 14:35:49.056636",
  "2021-10-05 14:35:49.056639",
  "2021-10-05 14:35:49.056641",
  "2021-10-05 14:35:49.056643",
  "2021-10-05 14:35:49.056646",
  "2021-10-05 14:35:49.056648",
  "2021-10-05 14:35:49.056651",
  "2021-10-05 14:35:49.056653",
  "2021-10-05 14:35:49.056656",
  "2021-10-05 14:35:49.056659",
  "2021-10-05 14:35:49.056662",
  "2021-10-05 14:35:49.056664",
  "2021-10-05 14:35:49.056667",
  "2021-10-05 14:35:49.056669",
  "2021-10-05 14:35:49.056672",
  "2021-10-05 14:35:49.056674",
  "2021-10-05 14:35:49.056677",
  "2021-10-05 14:35:49.056679",
  "2021-10-05 14:35:49.056681",
  "2021-10-05 14:35:49.056684",
  "2021-10-05 14:35:49.056686",
  "2021-10-05 14:35:49.056691",
  "2021-10-05 14:35:49.056694",
  "2021-10-05 14:35:49.056697",
  "2021-10-05 14:35:49.056699",
  "2021-10-05 14:35:49.056701",
  "2021-10-05 14:35:49.056704",
  "2021-10-05 14:35:49.056706",
  "2021-10-05 14:35:49.056709",
  "2021-10-05 14:35:49.056711",
  "2021-10-05 14:35:49.056714",
  "2021-10-05 14:35:49.056716",
  "2021-10-05 14:35:49.056719",
  "2021-10-05 14:35:49.056721",
  "2021-10-05 14:35:49.056723",
  "2021-10-05 14:35:49.056726",
  "2021-10-05 14:35:49.056728",
  "2021-10-05 14:35:49.056732",
  "2021-10-05 14:35:49.056735",
  "2021-10-05 14:35:49.056738",
  "2021-10-05 14:35:49.056740",
  "2021-10-05 14:35:49.056742",
  "2021-10-05 14:35:49.056745",
  "2021-10-05 14:35:49.056747",
  "2021-10-05 14:35:49.056750",
  "2021-10-05 14:35:49.056752",
  "2021-10-05 14:35:49.056755",
  "2021-10-05 14:35:49.056757",
  "2021-10-05 14:35:49.056760",
  "2021-10-05 14:35:49.056762",
  "2021-10-05 14:35:49.056764",
  "2021-10-05 14:35:49.056767",
  "2021-10-05 14:35:49.056769",
  "2021-10-05 14:35:49.056774",
  "2021-10-05 14:35:49.056777",
  "2021-10-05 14:35:49.056779",
  "2021-10-05 14:35:49.056782",
  "2021-10-05 14:35:49.056784",
  "2021-10-05 14:35:49.056787",
  "2021-10-05 14:35:49.056789",
  "2021-10-05 14:35:49.056791",
)
$arr = New-Object "object[,]" 116,1
for ($i = 0; $i -lt 116; $i++) { $arr[$i,0] = $ts[$i] }
$ws1.Range("F2:F117").Value = $arr

# ---------------------------------------------------------------------
# 2. Add the new "metadata" sheet right after "data".
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Copy the header-row / index-column formatting from the "data" sheet so
# the new sheet re-uses the exact same bold+bordered style.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Syndromic Retinopathy"
$ws2.Range("C2").Value = 3099

# "0.178" must stay text (matches the source inlineStr), not become 0.178
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.178"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("E2").Value = "2021-07-22T23:45:18.524031Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:49.053054"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3099/?format=json"

# Match page-margin defaults of the rest of the workbook (1in top/bottom,
# 0.75in left/right, 0.5in header/footer).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36
